# Case_4_205/res_bus/vm_pu.xlsx -- "case with 380 kV done"
#
# The bus-voltage-magnitude (p.u.) results table occupies B2:N25
# (column A is the row/time index 0..23, column G stays fixed at 1,
# column H is intentionally blank). This rerun with the 380 kV slack
# voltage (B column 1.05 -> 1.02) shifts every other bus voltage in
# the table, so every data cell in B:F and I:N is rewritten below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.075127692472976
$ws.Range("D2").Value = 1.068859061594275
$ws.Range("E2").Value = 1.088373604859081
$ws.Range("F2").Value = 1.09453421040821
$ws.Range("I2").Value = 1.05657694330282
$ws.Range("J2").Value = 1.080033712190591
$ws.Range("K2").Value = 1.071563117270248
$ws.Range("L2").Value = 1.091026459752524
$ws.Range("M2").Value = 1.097171310500737
$ws.Range("N2").Value = 1.0815674839435

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.076498757581475
$ws.Range("D3").Value = 1.069899994862129
$ws.Range("E3").Value = 1.089731779934121
$ws.Range("F3").Value = 1.095953138768919
$ws.Range("I3").Value = 1.057029543662617
$ws.Range("J3").Value = 1.081062139016232
$ws.Range("K3").Value = 1.072419551209251
$ws.Range("L3").Value = 1.092203082345252
$ws.Range("M3").Value = 1.098409689955691
$ws.Range("N3").Value = 1.082597371253196

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.077385117009685
$ws.Range("D4").Value = 1.07057272828931
$ws.Range("E4").Value = 1.090610145486834
$ws.Range("F4").Value = 1.096870926463714
$ws.Range("I4").Value = 1.057320741832275
$ws.Range("J4").Value = 1.081726279981585
$ws.Range("K4").Value = 1.072972295338066
$ws.Range("L4").Value = 1.092963423344583
$ws.Range("M4").Value = 1.099210112239794
$ws.Range("N4").Value = 1.083262455374898

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.077757552791121
$ws.Range("D5").Value = 1.070855351786996
$ws.Range("E5").Value = 1.090979303111295
$ws.Range("F5").Value = 1.097256683877657
$ws.Range("I5").Value = 1.057442764624727
$ws.Range("J5").Value = 1.082005172091881
$ws.Range("K5").Value = 1.073204330216882
$ws.Range("L5").Value = 1.093282831917034
$ws.Range("M5").Value = 1.099546399999443
$ws.Range("N5").Value = 1.083541743543979

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.077820075403443
$ws.Range("D6").Value = 1.070902794238544
$ws.Range("E6").Value = 1.091041280126193
$ws.Range("F6").Value = 1.097321449643522
$ws.Range("I6").Value = 1.057463229558659
$ws.Range("J6").Value = 1.082051981058719
$ws.Range("K6").Value = 1.07324327010507
$ws.Range("L6").Value = 1.093336448124129
$ws.Range("M6").Value = 1.099602852002896
$ws.Range("N6").Value = 1.083588618984919

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.07739009425467
$ws.Range("D7").Value = 1.070576505475848
$ws.Range("E7").Value = 1.090615078605585
$ws.Range("F7").Value = 1.096876081284396
$ws.Range("I7").Value = 1.057322373864463
$ws.Range("J7").Value = 1.081730007776474
$ws.Range("K7").Value = 1.0729753971273
$ws.Range("L7").Value = 1.092967692231327
$ws.Range("M7").Value = 1.099214606552834
$ws.Range("N7").Value = 1.083266188463682

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.075591221073139
$ws.Range("D8").Value = 1.069211021266379
$ws.Range("E8").Value = 1.088832705056568
$ws.Range("F8").Value = 1.095013819672455
$ws.Range("I8").Value = 1.056730247313564
$ws.Range("J8").Value = 1.08038154911702
$ws.Range("K8").Value = 1.071852849945231
$ws.Range("L8").Value = 1.09142431668726
$ws.Range("M8").Value = 1.0975900136187
$ws.Range("N8").Value = 1.081915814838261

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.072414959582988
$ws.Range("D9").Value = 1.066798463001617
$ws.Range("E9").Value = 1.085688199619101
$ws.Range("F9").Value = 1.091729379372889
$ws.Range("I9").Value = 1.05567402564602
$ws.Range("J9").Value = 1.07799513738987
$ws.Range("K9").Value = 1.069863740374791
$ws.Range("L9").Value = 1.088696775153828
$ws.Range("M9").Value = 1.094720271441916
$ws.Range("N9").Value = 1.079526014132734

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.070292851921813
$ws.Range("D10").Value = 1.065185616854371
$ws.Range("E10").Value = 1.083589100593731
$ws.Range("F10").Value = 1.089537540343751
$ws.Range("I10").Value = 1.054961158984063
$ws.Range("J10").Value = 1.076397095226221
$ws.Range("K10").Value = 1.06853008011757
$ws.Range("L10").Value = 1.086872862820805
$ws.Range("M10").Value = 1.092802162743129
$ws.Range("N10").Value = 1.077925702565915

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.069372800741864
$ws.Range("D11").Value = 1.064486136859294
$ws.Range("E11").Value = 1.082679456037654
$ws.Range("F11").Value = 1.088587865758225
$ws.Range("I11").Value = 1.054650389192339
$ws.Range("J11").Value = 1.075703394657807
$ws.Range("K11").Value = 1.067950753438492
$ws.Range("L11").Value = 1.086081720473937
$ws.Range("M11").Value = 1.091970373160463
$ws.Range("N11").Value = 1.077231016863125

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.069030872203519
$ws.Range("D12").Value = 1.064226149115673
$ws.Range("E12").Value = 1.082341460308267
$ws.Range("F12").Value = 1.088235019916171
$ws.Range("I12").Value = 1.054534638865354
$ws.Range("J12").Value = 1.075445458215577
$ws.Range("K12").Value = 1.067735285514212
$ws.Range("L12").Value = 1.085787643689582
$ws.Range("M12").Value = 1.091661219067853
$ws.Range("N12").Value = 1.076972714121561

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.069104225296899
$ws.Range("D13").Value = 1.064281925104956
$ws.Range("E13").Value = 1.082413966743546
$ws.Range("F13").Value = 1.08831071090342
$ws.Range("I13").Value = 1.054559482062231
$ws.Range("J13").Value = 1.075500798523067
$ws.Range("K13").Value = 1.067781516848481
$ws.Range("L13").Value = 1.085850733792323
$ws.Range("M13").Value = 1.091727542362695
$ws.Range("N13").Value = 1.077028133018636

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.069344540538193
$ws.Range("D14").Value = 1.064464649673577
$ws.Range("E14").Value = 1.082651519543957
$ws.Range("F14").Value = 1.088558701367104
$ws.Range("I14").Value = 1.054640827704225
$ws.Range("J14").Value = 1.075682078983996
$ws.Range("K14").Value = 1.067932948527903
$ws.Range("L14").Value = 1.086057416331719
$ws.Range("M14").Value = 1.091944822298414
$ws.Range("N14").Value = 1.077209670918612

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.069492582595087
$ws.Range("D15").Value = 1.06457720971758
$ws.Range("E15").Value = 1.082797868510895
$ws.Range("F15").Value = 1.088711483809306
$ws.Range("I15").Value = 1.054690905427363
$ws.Range("J15").Value = 1.075793736599621
$ws.Range("K15").Value = 1.068026213354676
$ws.Range("L15").Value = 1.086184732131829
$ws.Range("M15").Value = 1.092078670265306
$ws.Range("N15").Value = 1.077321487100859

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.070353887674973
$ws.Range("D16").Value = 1.065232015423989
$ws.Range("E16").Value = 1.083649455098746
$ws.Range("F16").Value = 1.089600554163786
$ws.Range("I16").Value = 1.054981739465352
$ws.Range("J16").Value = 1.07644309688743
$ws.Range("K16").Value = 1.068568489009296
$ws.Range("L16").Value = 1.08692533895092
$ws.Range("M16").Value = 1.092857339436798
$ws.Range("N16").Value = 1.077971769554759

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.070893845863157
$ws.Range("D17").Value = 1.065642458968587
$ws.Range("E17").Value = 1.084183436398682
$ws.Range("F17").Value = 1.090158081860708
$ws.Range("I17").Value = 1.055163609977436
$ws.Range("J17").Value = 1.076849955313528
$ws.Range("K17").Value = 1.068908148447045
$ws.Range("L17").Value = 1.087389530389636
$ws.Range("M17").Value = 1.093345443827602
$ws.Range("N17").Value = 1.07837920576649

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.071208682356644
$ws.Range("D18").Value = 1.065881757170369
$ws.Range("E18").Value = 1.084494829443595
$ws.Range("F18").Value = 1.090483221392172
$ws.Range("I18").Value = 1.055269490083759
$ws.Range("J18").Value = 1.07708710155853
$ws.Range("K18").Value = 1.069106088454468
$ws.Range("L18").Value = 1.087660152919695
$ws.Range("M18").Value = 1.093630027940052
$ws.Range("N18").Value = 1.07861668878636

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.071316014639393
$ws.Range("D19").Value = 1.065963333699183
$ws.Range("E19").Value = 1.084600994833811
$ws.Range("F19").Value = 1.090594076068327
$ws.Range("I19").Value = 1.055305558281039
$ws.Range("J19").Value = 1.077167934030077
$ws.Range("K19").Value = 1.069173550858052
$ws.Range("L19").Value = 1.087752405873659
$ws.Range("M19").Value = 1.093727043739508
$ws.Range("N19").Value = 1.078697636049289

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.070835925084639
$ws.Range("D20").Value = 1.065598433309947
$ws.Range("E20").Value = 1.084126152455025
$ws.Range("F20").Value = 1.090098270344117
$ws.Range("I20").Value = 1.055144117889991
$ws.Range("J20").Value = 1.076806320599746
$ws.Range("K20").Value = 1.068871724598492
$ws.Range("L20").Value = 1.087339740778701
$ws.Range("M20").Value = 1.09329308718201
$ws.Range("N20").Value = 1.078335509086409

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.06927377871894
$ws.Range("D21").Value = 1.064410846555326
$ws.Range("E21").Value = 1.082581569330873
$ws.Range("F21").Value = 1.088485677004821
$ws.Range("I21").Value = 1.054616882195802
$ws.Range("J21").Value = 1.075628703759285
$ws.Range("K21").Value = 1.067888363434271
$ws.Range("L21").Value = 1.085996559344396
$ws.Range("M21").Value = 1.091880844067174
$ws.Range("N21").Value = 1.07715621989496

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.068290547810686
$ws.Range("D22").Value = 1.06366317991909
$ws.Range("E22").Value = 1.08160977098204
$ws.Range("F22").Value = 1.087471226707691
$ws.Range("I22").Value = 1.054283555495123
$ws.Range("J22").Value = 1.074886753116549
$ws.Range("K22").Value = 1.067268462429615
$ws.Range("L22").Value = 1.085150823144563
$ws.Range("M22").Value = 1.09099180677679
$ws.Range("N22").Value = 1.076413215597209

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.068811878285125
$ws.Range("D23").Value = 1.064059626539835
$ws.Range("E23").Value = 1.08212500366627
$ws.Range("F23").Value = 1.088009059772704
$ws.Range("I23").Value = 1.054460432679504
$ws.Range("J23").Value = 1.075280222172581
$ws.Range("K23").Value = 1.067597238615188
$ws.Range("L23").Value = 1.085599281510516
$ws.Range("M23").Value = 1.091463208496839
$ws.Range("N23").Value = 1.076807243424432

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.070862097336961
$ws.Range("D24").Value = 1.065618326939178
$ws.Range("E24").Value = 1.084152036817008
$ws.Range("F24").Value = 1.090125296769576
$ws.Range("I24").Value = 1.055152926148974
$ws.Range("J24").Value = 1.076826037765663
$ws.Range("K24").Value = 1.068888183516383
$ws.Range("L24").Value = 1.087362238969823
$ws.Range("M24").Value = 1.093316745263681
$ws.Range("N24").Value = 1.078355254252962

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.073236889090056
$ws.Range("D25").Value = 1.067422943145514
$ws.Range("E25").Value = 1.08650159804454
$ws.Range("F25").Value = 1.092578856391844
$ws.Range("I25").Value = 1.055948613660119
$ws.Range("J25").Value = 1.078613318467201
$ws.Range("K25").Value = 1.070379297962014
$ws.Range("L25").Value = 1.089402871650172
$ws.Range("M25").Value = 1.095463022936175
$ws.Range("N25").Value = 1.080145073098101

Write-Output "vm_pu.xlsx: updated B2:N25 for the 380 kV case"
